# "combobox de producto completado"
# Adds the newly-entered Product row (row 2) to the "Product" worksheet,
# mirroring a row picked from the Category/Brand combobox-style lookups.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "ssdgg"
$ws.Range("C2").Value = "Categoría 1"
$ws.Range("D2").Value = "Marca A"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = "fgj"
